$wb = $excel.ActiveWorkbook

# --- Update scraped_at timestamps on the "snapshot" sheet (column K) ---
$snapshot = $wb.Worksheets.Item("snapshot")

$snapshot.Range("K2").Value = "2025-12-17T03:58:00.087830+00:00"
$snapshot.Range("K3").Value = "2025-12-17T03:58:00.087851+00:00"
$snapshot.Range("K4").Value = "2025-12-17T03:58:02.595780+00:00"
$snapshot.Range("K5").Value = "2025-12-17T03:58:02.595822+00:00"
$snapshot.Range("K6").Value = "2025-12-17T03:58:02.595838+00:00"
$snapshot.Range("K7").Value = "2025-12-17T03:58:04.720958+00:00"
$snapshot.Range("K8").Value = "2025-12-17T03:58:07.527469+00:00"
$snapshot.Range("K9").Value = "2025-12-17T03:58:09.827258+00:00"
$snapshot.Range("K10").Value = "2025-12-17T03:58:09.827287+00:00"
$snapshot.Range("K11").Value = "2025-12-17T03:58:12.092570+00:00"
$snapshot.Range("K12").Value = "2025-12-17T03:58:16.722914+00:00"
$snapshot.Range("K13").Value = "2025-12-17T03:58:16.722942+00:00"
$snapshot.Range("K14").Value = "2025-12-17T03:58:19.043266+00:00"
$snapshot.Range("K15").Value = "2025-12-17T03:58:21.786863+00:00"
$snapshot.Range("K16").Value = "2025-12-17T03:58:24.058827+00:00"
$snapshot.Range("K17").Value = "2025-12-17T03:58:26.352621+00:00"
$snapshot.Range("K18").Value = "2025-12-17T03:58:26.352651+00:00"
$snapshot.Range("K19").Value = "2025-12-17T03:58:26.352673+00:00"
$snapshot.Range("K20").Value = "2025-12-17T03:58:26.352691+00:00"
$snapshot.Range("K21").Value = "2025-12-17T03:58:28.664492+00:00"
$snapshot.Range("K22").Value = "2025-12-17T03:58:28.664520+00:00"
$snapshot.Range("K23").Value = "2025-12-17T03:58:30.919249+00:00"
$snapshot.Range("K24").Value = "2025-12-17T03:58:30.919276+00:00"
$snapshot.Range("K25").Value = "2025-12-17T03:58:30.919293+00:00"
$snapshot.Range("K26").Value = "2025-12-17T03:58:33.676753+00:00"
$snapshot.Range("K27").Value = "2025-12-17T03:58:36.474339+00:00"
$snapshot.Range("K28").Value = "2025-12-17T03:58:36.474369+00:00"
$snapshot.Range("K29").Value = "2025-12-17T03:58:36.474386+00:00"
$snapshot.Range("K30").Value = "2025-12-17T03:58:38.778458+00:00"
$snapshot.Range("K31").Value = "2025-12-17T03:58:41.623372+00:00"
$snapshot.Range("K32").Value = "2025-12-17T03:58:41.623399+00:00"
$snapshot.Range("K33").Value = "2025-12-17T03:58:46.252284+00:00"
$snapshot.Range("K34").Value = "2025-12-17T03:58:46.252315+00:00"
$snapshot.Range("K35").Value = "2025-12-17T03:58:48.980678+00:00"
$snapshot.Range("K36").Value = "2025-12-17T03:58:48.980705+00:00"

# --- Remove the now-stale rows from the "new_injured" sheet, keeping only the header ---
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Range("A2:G4").Delete()

Write-Output "done"
